# Sync attendance_reports: normalise "Recorded By" (column G) ordering.
# For every row whose G cell lists multiple comma-separated recorders and
# does NOT already start with "backup@backdoor.com", rotate the first
# name to the end of the list (matches the upstream canonical ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dim = $ws.UsedRange
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 157 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    if ([string]::IsNullOrEmpty($current)) { continue }

    $parts = $current -split ',\s*'
    if ($parts.Count -gt 1 -and $parts[0] -ne 'backup@backdoor.com') {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ', '
        $cell.Value = $rotated
    }
}

Write-Host "Recorded By column normalised"
